# Append the daily rows that follow 29-10-2021 in the "Recompra deuda BCCh en
# pesos 2021 - Diaria" sheet: 30-10-2021, 31-10-2021, 01-11-2021, 02-11-2021,
# each carrying the same BCP/BCU/Otros figures (696 / 3962 / 59) as the rows
# before them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDates = @("30-10-2021", "31-10-2021", "01-11-2021", "02-11-2021")

$row = 304
foreach ($d in $newDates) {
    $cell = $ws.Range("A$row")

    # Some of these strings (e.g. "01-11-2021") are ambiguous enough that
    # Excel's smart text-to-date parsing would otherwise store them as real
    # date serials instead of literal text like the rest of column A. Enter
    # the text via a formula (never auto-converted) and then collapse it to
    # a plain value with copy / paste-values, which keeps it literal text
    # without leaving behind any new number-format/style definition.
    $cell.Formula = '="' + $d + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
    $excel.CutCopyMode = $false

    $ws.Range("B$row").Value = 696
    $ws.Range("C$row").Value = 3962
    $ws.Range("D$row").Value = 59

    $row = $row + 1
}
